$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data values in row 4
$ws.Range("E4").Value = 5
$ws.Range("G4").Value = -3
$ws.Range("H4").Value = 13

# Update the active selection to match the saved view state
$ws.Range("E4").Select()
